$wb = $excel.ActiveWorkbook

# Table_Names sheet: A1 "T2" -> "T1"
$wsTable = $wb.Worksheets.Item("Table_Names")
$wsTable.Range("A1").Value = "T1"

# Field_Names sheet: A1 (empty) -> "T1F2"
$wsField = $wb.Worksheets.Item("Field_Names")
$wsField.Range("A1").Value = "T1F2"
